$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-8, columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$rows = @{
    2 = @{ D = 44179; M = 45; N = 20000; O = 20000; P = 20000; S = 4000 }
    3 = @{ D = 44175; M = 25; N = 20000; O = 20000; P = 20000; S = 4000 }
    4 = @{ D = 44186; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    5 = @{ D = 44188; M = 30; N = 15000; O = 15000; P = 15000; S = 3000 }
    6 = @{ D = 44193; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    7 = @{ D = 44181; M = 30; N = 20000; O = 20000; P = 20000; S = 4000 }
    8 = @{ D = 44196; M = 56; N = 15000; O = 15000; P = 15000; S = 3000 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("S$r").Value = $vals.S
}
